$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45171 = 2023-09-02)
# that was bumped by one day (45172 = 2023-09-03) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 116 }

$ws.Range("C2:C$lastRow").Value = 45172
